$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.744.09"
$ws.Range("E2").Value = "  +2.82%  "
$ws.Range("D3").Value = "2.419.65"
$ws.Range("E3").Value = "  +2.57%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'551.85"
$ws.Range("E5").Value = "  +0.10%  "
$ws.Range("D6").Value = "'137.05"
$ws.Range("E6").Value = "  +3.30%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E8").Value = "  +3.71%  "
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("D10").Value = "'5.68"
$ws.Range("E10").Value = "  +0.87%  "
$ws.Range("E11").Value = "  -2.19%  "
$ws.Range("D12").Value = "'0.354"
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  +3.20%  "
$ws.Range("D14").Value = "2.849.71"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").Value = "59.707.29"
$ws.Range("E15").Value = "  +3.02%  "
$ws.Range("E16").Value = "  +0.87%  "
$ws.Range("D17").Value = "2.434.85"
$ws.Range("E17").Value = "  +3.06%  "
$ws.Range("E18").Value = "  +2.39%  "
$ws.Range("E19").Value = "  +1.76%  "
$ws.Range("D20").Value = "'329.76"
$ws.Range("E20").Value = "  -0.30%  "
$ws.Range("E21").Value = "  -3.11%  "
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'65.98"
$ws.Range("E23").Value = "  +3.21%  "
$ws.Range("E24").Value = "  +1.60%  "
$ws.Range("D25").Value = "'8.76"
$ws.Range("E25").Value = "  +5.89%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("E27").Value = "  +3.15%  "
$ws.Range("D28").Value = "0.0₃0773"
$ws.Range("E28").Value = "  +4.77%  "
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").Value = "'170.51"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("E31").Value = "  -0.19%  "
$ws.Range("E32").Value = "  +1.42%  "
$ws.Range("E33").Value = "  +0.40%  "
$ws.Range("D34").Value = "'0.998"
$ws.Range("E35").Value = "  +4.05%  "
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "'39.58"
$ws.Range("E39").Value = "  -2.01%  "
$ws.Range("D40").Value = "'0.410"
$ws.Range("E40").Value = "  -5.50%  "
$ws.Range("D41").Value = "'312.76"
$ws.Range("E41").Value = "  +7.86%  "
$ws.Range("E42").Value = "  +0.21%  "
$ws.Range("D43").Value = "'138.63"
$ws.Range("E43").Value = "  -2.31%  "
$ws.Range("D44").Value = "'0.0971"
$ws.Range("E44").Value = "  +1.69%  "
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").Value = "'19.51"
$ws.Range("E46").Value = "  +3.90%  "
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("B48").Value = "Polygon"
$ws.Range("C48").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D48").Value = "'0.404"
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("B49").Value = "VeChain"
$ws.Range("C49").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D49").Value = "'0.0224"
$ws.Range("E49").Value = "  +1.04%  "
$ws.Range("E50").Value = "  +0.73%  "
$ws.Range("D51").Value = "'11.03"
$ws.Range("E51").Value = "  -0.27%  "
